# Applies the re-ordering of the news-article rows (A2:E6) so that each
# row's title/timestamp/historical-distance/time-bucket/uri once again line
# up with one another after a new JSON source was merged in for the time
# bucket analysis. Row 6 (Epping murder) keeps its original values; rows
# 2-5 are rewritten with the same four articles but in their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title2 = "Man accused of 2009 Lin family murders, Robert Xie, granted bail in NSW Supreme Court"
$time2  = "2015-12-08T05:29:25UTC"
$dist2  = 2334
$bucket2 = "day_31_beyond"
$uri2   = "http://www.smh.com.au/nsw/man-accused-of-2009-lin-family-murders-robert-xie-granted-bail-in-nsw-supreme-court-20151208-gli4d7.html"

$title3 = "Police video captured Robert Xie destroying evidence in Lin family case"
$time3  = "2017-01-16T18:14:58UTC"
$dist3  = 2739
$bucket3 = "day_31_beyond"
$uri3   = "http://www.smh.com.au/nsw/police-video-captured-robert-xie-destroying-evidence-in-lin-family-case-20170116-gtsfx0"

$title4 = "Robert Xie murder trial: the evidence that helped catch a killer"
$time4  = "2017-01-16T05:25:20UTC"
$dist4  = 2739
$bucket4 = "day_31_beyond"
$uri4   = "http://www.abc.net.au/news/2017-01-16/robert-xie-trial-blood-stains-helped-convict-murderer/8184856"

$title5 = "A motive in Lin murders revealed: Robert Xie's niece Brenda Lin breaks her silence"
$time5  = "2017-02-16T14:06:37UTC"
$dist5  = 2770
$bucket5 = "day_31_beyond"
$uri5   = "http://www.smh.com.au/nsw/motive-in-lin-murders-revealed:-robert-xie's-niece-brenda-lin-breaks-her-silence-20170216-gue88t.html"

$title6 = "Epping murder"
$time6  = "2009-07-19T00:00:00UTC"
$dist6  = 1
$bucket6 = "day_1"
$uri6   = "http://www.smh.com.au/national/grisly-killings-that-shocked-even-the-detectives-20090719-dpku.html"

$ws.Range("A2").Value = $title2
$ws.Range("B2").Value = $time2
$ws.Range("C2").Value = $dist2
$ws.Range("D2").Value = $bucket2
$ws.Range("E2").Value = $uri2

$ws.Range("A3").Value = $title3
$ws.Range("B3").Value = $time3
$ws.Range("C3").Value = $dist3
$ws.Range("D3").Value = $bucket3
$ws.Range("E3").Value = $uri3

$ws.Range("A4").Value = $title4
$ws.Range("B4").Value = $time4
$ws.Range("C4").Value = $dist4
$ws.Range("D4").Value = $bucket4
$ws.Range("E4").Value = $uri4

$ws.Range("A5").Value = $title5
$ws.Range("B5").Value = $time5
$ws.Range("C5").Value = $dist5
$ws.Range("D5").Value = $bucket5
$ws.Range("E5").Value = $uri5

$ws.Range("A6").Value = $title6
$ws.Range("B6").Value = $time6
$ws.Range("C6").Value = $dist6
$ws.Range("D6").Value = $bucket6
$ws.Range("E6").Value = $uri6

# The hyperlinks attached to column E need to follow the article that now
# sits in each row. The COM layer only supports editing hyperlinks by
# recreating them, so drop the whole collection and rebuild it against the
# correct address for each row, then restore the Hyperlink cell style.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), $uri2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), $uri3) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), $uri4) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), $uri5) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"), $uri6) | Out-Null

$ws.Range("E2:E6").Style = "Hyperlink"
